$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), values for columns B-E and G.
# Column F ("Win") is left unchanged.
$data = @{
    2  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 0.496779210170732; G = 9.295990156953671 }
    3  = @{ B = 0.3048080303191223; C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 3.274871460341982 }
    4  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732; G = 5.553084769722144 }
    5  = @{ B = 0.01514828764759746;C = 0.002777888934908601;D = 0.8054896365839992; E = 0.496779210170732; G = 1.320195023337237 }
    6  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    8  = @{ B = 0.01514828764759746;C = 0.002777888934908601;D = 0.8054896365839992; E = 0.496779210170732; G = 1.320195023337237 }
    9  = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;   E = 8.660232485948974; G = 17.45944343273191 }
    10 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
